$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -8.016
$ws.Range("A9").Value = -21.847
$ws.Range("D9").Value = -8.123000000000001
$ws.Range("D11").Value = -7.323
$ws.Range("A18").Value = -22.095
$ws.Range("A20").Value = -20.631
$ws.Range("D23").Value = -8.387
$ws.Range("D24").Value = -7.101000000000001
$ws.Range("D26").Value = -7.48
$ws.Range("A27").Value = -21.922
$ws.Range("D34").Value = -7.606999999999999
$ws.Range("A35").Value = -20.089
$ws.Range("D35").Value = -8.254000000000001
$ws.Range("D48").Value = -7.892999999999999
$ws.Range("D49").Value = -8.35
$ws.Range("D52").Value = -7.958000000000001
$ws.Range("D66").Value = -7.486
$ws.Range("D67").Value = -7.616
$ws.Range("A69").Value = -21.52
$ws.Range("A76").Value = -20.074
$ws.Range("A78").Value = -19.911
$ws.Range("D78").Value = -8.379999999999999
$ws.Range("D80").Value = -8.068999999999999
$ws.Range("A82").Value = -21.901
$ws.Range("A83").Value = -21.891
$ws.Range("A93").Value = -21.461
$ws.Range("D99").Value = -8.253
$ws.Range("D104").Value = -7.647
